$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Uni")

# Rotate fuel codes in column A across the repeated ~TradeLinks blocks
$ws.Range("A2").Value  = "FOL"
$ws.Range("A7").Value  = "JET"
$ws.Range("A12").Value = "OTH"
$ws.Range("A17").Value = "COA"
$ws.Range("A22").Value = "PET"
$ws.Range("A27").Value = "COL"
$ws.Range("A32").Value = "DID"
$ws.Range("A37").Value = "DIJ"
$ws.Range("A42").Value = "LPG"
$ws.Range("A47").Value = "DSL"

# Move the stray "1" marker from B34 to B19
$ws.Range("B34").ClearContents()
$ws.Range("B19").Value = "1"
